# Revert "Drop in all data files from 3.0 RMI script"
#
# Re-introduces a "Texas Notes" worksheet between "About" and "DR", restores
# the discount-rate value to 5.87% (from the placeholder 3%), and adds the
# two supporting notes explaining the change.

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$drSheet    = $wb.Worksheets.Item("DR")

# Insert a brand-new sheet immediately before "DR" and rename it.
$texasSheet = $wb.Worksheets.Add($drSheet)
$texasSheet.Name = "Texas Notes"

# Re-resolve sheet handles by name: inserting/renaming a sheet can leave
# older variable handles pointing at the wrong sheet object.
$aboutSheet = $wb.Worksheets.Item("About")
$drSheet    = $wb.Worksheets.Item("DR")
$texasSheet = $wb.Worksheets.Item("Texas Notes")

# Populate the new "Texas Notes" sheet.
$texasSheet.Range("A1").Value = "updated to the VCE WISdom number"
$texasSheet.Range("A2").Value = 0.0587
$texasSheet.Range("A4").Value = "their feedback was the 3% was a bit low"

# Restore the discount rate on the "DR" sheet.
$drSheet.Range("B2").Value = 0.0587

# Restore selections on each sheet.
$aboutSheet.Activate()
$aboutSheet.Range("C23").Select() | Out-Null

$texasSheet.Activate()
$texasSheet.Range("A5").Select() | Out-Null

$drSheet.Activate()
$drSheet.Range("B2").Select() | Out-Null

# "About" is the tab that is actually selected/visible.
$aboutSheet.Activate()

# Turn on iterative calculation (calcPr iterate="1" iterateDelta="1e-5").
$excel.Iteration = $true
$excel.MaxChange = 0.00001
